$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Process column (H2) changes from "Search" to "Search & Typing"
$ws.Range("H2").Value = "Search & Typing"

# Row 3: Process column (H3) changes from "Typing" to "Search & Typing"
$ws.Range("H3").Value = "Search & Typing"

# Update the active selection to K5 (matches recorded sheetView selection)
$ws.Range("K5").Select()
